$d = $word.ActiveDocument

# --- Change 1: remove the "Fecha del informe..." paragraph and the
# following space-only paragraph entirely. The paragraph that follows
# keeps its own (unrelated) formatting untouched. ---
$dateLabel = "Fecha del informe: 22 de enero de 2024"
$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match [regex]::Escape($dateLabel)) {
        $startPara = $p
        break
    }
}

if ($startPara -ne $null) {
    $nextPara = $startPara.Next()
    $rangeToDelete = $d.Range($startPara.Range.Start, $nextPara.Range.End)
    $rangeToDelete.Delete()
}

# --- Change 2: reword the closing sentence of the report. The replacement
# is written directly via InsertXML on a freshly-built Range(start,end)
# (with the original run's formatting preserved) so that neighboring runs
# in the paragraph are not coalesced the way a plain Find/Replace would
# coalesce them. ---
$oldSentence = "Si el producto puede continuar su ritmo de ventas actual en la meca del mundo de la salud y el ejercicio, entonces podría estar listo para la exposición a nivel nacional."
$newSentence = "Si el producto puede continuar el ritmo de ventas actual en la meca del mundo de la salud y el ejercicio, puede estar listo para la exposición a nivel nacional."

$hit = $d.Content
$found = $hit.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # Rebuild a plain Range over the exact same span the Find located -
    # using the Find-owned range object directly for InsertXML inserts
    # rather than replaces, so we re-anchor it first.
    $targetRange = $d.Range($hit.Start, $hit.End)

    $runProps = '<w:rPr><w:rStyle w:val="DefaultParagraphFont"/><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:dstrike w:val="0"/><w:outline w:val="0"/><w:shadow w:val="0"/><w:emboss w:val="0"/><w:imprint w:val="0"/><w:noProof w:val="0"/><w:vanish w:val="0"/><w:color w:val="auto"/><w:spacing w:val="0"/><w:w w:val="100"/><w:kern w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:highlight w:val="none"/><w:u w:val="none" w:color="auto"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto"/><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vertAlign w:val="baseline"/><w:rtl w:val="0"/><w:cs w:val="0"/><w:lang w:val="es-ES" w:eastAsia="ja-JP" w:bidi="ar-SA"/></w:rPr>'

    $xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $runProps + '<w:t>' + $newSentence + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $targetRange.InsertXML($xmlFrag)
}
